# Apply updated cryptocurrency price/volume figures (and the two row
# re-orderings: TRON/WrappedEther swap rows 15-16, Litecoin/ICP swap rows
# 24-25, Kaspa/Dai swap rows 35-36) per the Fri Mar 22 2024 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.862.53'
$ws.Range('E2').Value = '  -3.03%  '

$ws.Range('D3').Value = '3.438.73'
$ws.Range('E3').Value = '  -2.81%  '

$ws.Range('D5').Value = '''572.05'
$ws.Range('E5').Value = '  +0.51%  '

$ws.Range('D6').Value = '''174.81'
$ws.Range('E6').Value = '  -7.40%  '

$ws.Range('D7').Value = '''0.624'
$ws.Range('E7').Value = '  +0.79%  '

$ws.Range('E8').Value = '  -0.07%  '

$ws.Range('D9').Value = '''0.625'
$ws.Range('E9').Value = '  -1.60%  '

$ws.Range('D10').Value = '''0.158'
$ws.Range('E10').Value = '  +4.49%  '

$ws.Range('D11').Value = '''55.24'
$ws.Range('E11').Value = '  +0.63%  '

$ws.Range('D12').Value = '''0.0000274'
$ws.Range('E12').Value = '  +0.95%  '

$ws.Range('D13').Value = '''9.13'
$ws.Range('E13').Value = '  -3.18%  '

$ws.Range('D14').Value = '4.002.55'
$ws.Range('E14').Value = '  -2.53%  '

$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '3.448.80'
$ws.Range('E15').Value = '  -2.66%  '

$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').Value = '''0.121'
$ws.Range('E16').Value = '  -0.71%  '

$ws.Range('D17').Value = '''18.13'
$ws.Range('E17').Value = '  -0.78%  '

$ws.Range('D18').Value = '''11.90'
$ws.Range('E18').Value = '  -1.05%  '

$ws.Range('D19').Value = '64.890.50'
$ws.Range('E19').Value = '  -3.03%  '

$ws.Range('D20').Value = '''0.992'
$ws.Range('E20').Value = '  -0.53%  '

$ws.Range('D21').Value = '''406.91'
$ws.Range('E21').Value = '  -5.54%  '

$ws.Range('D22').Value = '''4.21'
$ws.Range('E22').Value = '  -0.47%  '

$ws.Range('D23').Value = '''4.45'
$ws.Range('E23').Value = '  +8.00%  '

$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '''83.77'
$ws.Range('E24').Value = '  -1.62%  '

$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').Value = '''13.36'
$ws.Range('E25').Value = '  +8.68%  '

$ws.Range('D26').Value = '''10.84'
$ws.Range('E26').Value = '  -2.58%  '

$ws.Range('D27').Value = '''2.81'
$ws.Range('E27').Value = '  -3.36%  '

$ws.Range('D28').Value = '''9.03'
$ws.Range('E28').Value = '  -2.51%  '

$ws.Range('D29').Value = '''29.84'
$ws.Range('E29').Value = '  -1.61%  '

$ws.Range('D30').Value = '''6.59'
$ws.Range('E30').Value = '  -0.42%  '

$ws.Range('D31').Value = '''11.55'
$ws.Range('E31').Value = '  -1.89%  '

$ws.Range('D32').Value = '''586.80'
$ws.Range('E32').Value = '  -8.64%  '

$ws.Range('E33').Value = '  -3.48%  '

$ws.Range('D34').Value = '''59.66'
$ws.Range('E34').Value = '  -0.35%  '

$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '''0.153'
$ws.Range('E35').Value = '  +2.31%  '

$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D36').Value = '''1.00'
$ws.Range('E36').Value = '  +0.06%  '

$ws.Range('D37').Value = '0.0₃0778'
$ws.Range('E37').Value = '  -4.74%  '

$ws.Range('D38').Value = '''3.54'
$ws.Range('E38').Value = '  +5.57%  '

$ws.Range('D39').Value = '''36.22'
$ws.Range('E39').Value = '  -6.45%  '

$ws.Range('D40').Value = '''0.377'
$ws.Range('E40').Value = '  -3.97%  '

$ws.Range('D41').Value = '3.191.90'
$ws.Range('E41').Value = '  +4.45%  '

$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  +0.22%  '

$ws.Range('D43').Value = '''2.93'
$ws.Range('E43').Value = '  +1.18%  '

$ws.Range('D44').Value = '''2.51'
$ws.Range('E44').Value = '  -6.17%  '

$ws.Range('D45').Value = '''3.24'
$ws.Range('E45').Value = '  -3.40%  '

$ws.Range('D46').Value = '''0.0410'
$ws.Range('E46').Value = '  -2.66%  '

$ws.Range('E47').Value = '  -1.10%  '

$ws.Range('D48').Value = '''2.62'
$ws.Range('E48').Value = '  -5.47%  '

$ws.Range('D49').Value = '''8.47'
$ws.Range('E49').Value = '  -2.31%  '

$ws.Range('D50').Value = '''136.66'
$ws.Range('E50').Value = '  -4.73%  '

$ws.Range('D51').Value = '''2.33'
$ws.Range('E51').Value = '  -3.23%  '
